{"js": "// Update the worksheet date and each \"AAA\u00d7B=\" arithmetic prompt to its\n// new value. Every \"old\" string is unique in the document body, so a scoped\n// body.search(...) + Range.insertText(..., \"Replace\") round-trip per pair\n// reproduces the diff exactly (one hit each, run formatting untouched).\nconst replacements = [\n  [\"2024-07-12 Friday\", \"2024-07-13 Saturday\"],\n  [\"560\u00d72=\", \"822\u00d79=\"],\n  [\"565\u00d77=\", \"417\u00d78=\"],\n  [\"727\u00d76=\", \"761\u00d73=\"],\n  [\"408\u00d79=\", \"357\u00d78=\"],\n  [\"897\u00d78=\", \"666\u00d79=\"],\n  [\"812\u00d79=\", \"554\u00d79=\"],\n  [\"573\u00d77=\", \"902\u00d75=\"],\n  [\"227\u00d73=\", \"599\u00d77=\"],\n  [\"636\u00d76=\", \"977\u00d79=\"],\n  [\"769\u00d76=\", \"910\u00d76=\"],\n  [\"645\u00d79=\", \"836\u00d77=\"],\n  [\"655\u00d76=\", \"831\u00d73=\"],\n  [\"699\u00d78=\", \"537\u00d78=\"],\n  [\"976\u00d79=\", \"889\u00d72=\"],\n  [\"826\u00d74=\", \"925\u00d74=\"],\n  [\"811\u00d75=\", \"965\u00d72=\"],\n  [\"843\u00d78=\", \"870\u00d77=\"],\n  [\"163\u00d78=\", \"310\u00d72=\"],\n  [\"522\u00d73=\", \"310\u00d76=\"],\n  [\"762\u00d79=\", \"676\u00d79=\"],\n  [\"583\u00d75=\", \"521\u00d74=\"],\n  [\"370\u00d77=\", \"781\u00d76=\"],\n  [\"172\u00d72=\", \"596\u00d73=\"],\n  [\"797\u00d73=\", \"236\u00d72=\"],\n  [\"170\u00d77=\", \"403\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const result of results.items) {\n    result.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and each \"AAA\u00d7B=\" arithmetic prompt to its\n# new value. Every \"old\" string is unique in the document body, so a\n# Find/Replace pass per pair (wdReplaceAll, but each only ever matches once)\n# is safe and mirrors the canonical diff exactly.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-07-12 Friday', '2024-07-13 Saturday'),\n    @('560\u00d72=', '822\u00d79='),\n    @('565\u00d77=', '417\u00d78='),\n    @('727\u00d76=', '761\u00d73='),\n    @('408\u00d79=', '357\u00d78='),\n    @('897\u00d78=', '666\u00d79='),\n    @('812\u00d79=', '554\u00d79='),\n    @('573\u00d77=', '902\u00d75='),\n    @('227\u00d73=', '599\u00d77='),\n    @('636\u00d76=', '977\u00d79='),\n    @('769\u00d76=', '910\u00d76='),\n    @('645\u00d79=', '836\u00d77='),\n    @('655\u00d76=', '831\u00d73='),\n    @('699\u00d78=', '537\u00d78='),\n    @('976\u00d79=', '889\u00d72='),\n    @('826\u00d74=', '925\u00d74='),\n    @('811\u00d75=', '965\u00d72='),\n    @('843\u00d78=', '870\u00d77='),\n    @('163\u00d78=', '310\u00d72='),\n    @('522\u00d73=', '310\u00d76='),\n    @('762\u00d79=', '676\u00d79='),\n    @('583\u00d75=', '521\u00d74='),\n    @('370\u00d77=', '781\u00d76='),\n    @('172\u00d72=', '596\u00d73='),\n    @('797\u00d73=', '236\u00d72='),\n    @('170\u00d77=', '403\u00d78=')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
